$wb = $excel.ActiveWorkbook

$sDataCombined = $wb.Worksheets.Item(1)      # DataCombined
$sPlotConfig   = $wb.Worksheets.Item(2)      # plotConfiguration
$sPlotGrids    = $wb.Worksheets.Item(3)      # plotGrids
$sPlotTypes    = $wb.Worksheets.Item(5)      # plotTypes

# ---------------------------------------------------------------------------
# Populate new cell values.
# The order below is chosen so that new shared-string entries are created in
# the same order as in the target workbook:
#   P1, "Aciclovri observed", "Aciclovir_Laskin 1982...", "Aciclovir",
#   "P2", "P1, P2", "P3", "Aciclovr2", "blabla", "1;2;3"
# ---------------------------------------------------------------------------

# plotConfiguration!A2 -> "P1"
$sPlotConfig.Range("A2").Value = "P1"

# DataCombined row 3 (new observed data row)
$sDataCombined.Range("A3").Value = "AciclovirPVB"
$sDataCombined.Range("B3").Value = "observed"
$sDataCombined.Range("B3").Style = "Normal"
$sDataCombined.Range("C3").Value = "Aciclovri observed"
$sDataCombined.Range("F3").Value = "Aciclovir_Laskin 1982.Group A_1_Human_PeripheralVenousBlood_Plasma_2.5 mg/kg_iv"

# DataCombined row 4, only B4 gets a (later cleared) value so the row exists
$sDataCombined.Range("B4").Value = "x"
$sDataCombined.Range("B4").ClearContents() | Out-Null
$sDataCombined.Range("B4").Style = "Normal"

# plotGrids!A2:B2 -> "Aciclovir" / "P1, P2"
$sPlotGrids.Range("A2").Value = "Aciclovir"

# plotConfiguration!A3 -> "P2"
$sPlotConfig.Range("A3").Value = "P2"

# plotGrids!B2 -> "P1, P2"
$sPlotGrids.Range("B2").Value = "P1, P2"

# plotConfiguration!A4 -> "P3"
$sPlotConfig.Range("A4").Value = "P3"

# plotGrids row 3
$sPlotGrids.Range("A3").Value = "Aciclovr2"
$sPlotGrids.Range("B3").Value = "P2"

# plotConfiguration extra column N
$sPlotConfig.Range("N1").Value = "blabla"
$sPlotConfig.Range("N2").Value = "1;2;3"

# ---------------------------------------------------------------------------
# Fill in the remaining (reused-string) cells for the new rows.
# ---------------------------------------------------------------------------
$sPlotConfig.Range("B2").Value = "AciclovirPVB"
$sPlotConfig.Range("C2").Value = "individual"

$sPlotConfig.Range("B3").Value = "AciclovirPVB"
$sPlotConfig.Range("C3").Value = "observedVsSimulated"

$sPlotConfig.Range("B4").Value = "AciclovirPVB"
$sPlotConfig.Range("C4").Value = "residualsVsSimulated"

# ---------------------------------------------------------------------------
# Style fix-up: D2/E2 on DataCombined lose their "applyProtection" style.
# ---------------------------------------------------------------------------
$sDataCombined.Range("D2").Style = "Normal"
$sDataCombined.Range("E2").Style = "Normal"

# Approximate the column E autofit that Excel performed after the edit.
$sDataCombined.Columns("E:E").AutoFit() | Out-Null

# ---------------------------------------------------------------------------
# Sheet view / selection changes.
# ---------------------------------------------------------------------------
$sPlotGrids.Activate() | Out-Null
$sPlotGrids.Range("B2").Select() | Out-Null

$sPlotTypes.Activate() | Out-Null
$sPlotTypes.Range("A2").Select() | Out-Null

$sDataCombined.Activate() | Out-Null
$sDataCombined.Range("E7").Select() | Out-Null

# plotConfiguration becomes the active (selected) tab, matching activeTab="1"
$sPlotConfig.Activate() | Out-Null
$sPlotConfig.Range("R6").Select() | Out-Null
